# PercianBattleCard.xlsx: level up Percian from Level 1 -> Level 2 Striker.
# Updates defenses/skills/abilities text and the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: rewrite a cell that starts with a bold "Label" run followed by
# one or more plain-text runs, changing ONLY the plain-text portions while
# re-asserting bold on the label run(s) so the rich-text split survives the
# whole-string rewrite. $segments is an ordered list of @{Text=...; Bold=$true/$false}
function Set-RichCell($range, $segments) {
    $full = ""
    foreach ($seg in $segments) {
        $full = $full + $seg.Text
    }
    $range.Value = $full
    $pos = 1
    foreach ($seg in $segments) {
        $len = $seg.Text.Length
        if ($seg.Bold -and $len -gt 0) {
            $range.Characters($pos, $len).Font.Bold = $true
        }
        $pos += $len
    }
}

# D1: Level 1 Striker -> Level 2 Striker
$ws.Range("D1").Value = "Level 2 Striker"

# A3: Initiative +3 -> Initiative +4  (bold "Initiative" + plain " +4")
Set-RichCell $ws.Range("A3") @(
    @{Text="Initiative"; Bold=$true},
    @{Text=" +4"; Bold=$false}
)

# C3: Senses Perception +3; low-light vision -> ...+4...
Set-RichCell $ws.Range("C3") @(
    @{Text="Senses"; Bold=$true},
    @{Text=" Perception +4; low-light vision"; Bold=$false}
)

# A4: HP 23; Bloodied 11 -> HP 28; Bloodied 14
Set-RichCell $ws.Range("A4") @(
    @{Text="HP"; Bold=$true},
    @{Text=" 28; "; Bold=$false},
    @{Text="Bloodied"; Bold=$true},
    @{Text=" 14"; Bold=$false}
)

# A5: AC 13; Fortitude 11; Reflex 13; Will 15 -> AC 14; Fortitude 12; Reflex 14; Will 16
Set-RichCell $ws.Range("A5") @(
    @{Text="AC"; Bold=$true},
    @{Text=" 14; "; Bold=$false},
    @{Text="Fortitude"; Bold=$true},
    @{Text=" 12; "; Bold=$false},
    @{Text="Reflex"; Bold=$true},
    @{Text=" 14; "; Bold=$false},
    @{Text="Will"; Bold=$true},
    @{Text=" 16"; Bold=$false}
)

# A6 (Speed 7) is unchanged.

# Power text updates (plain strings; +1 modifier/values across the board)
$ws.Range("B8").Value = "Ranged 10; 4 vs Will; 1d10 + 4 psychic damage. If I roll an"
$ws.Range("B12").Value = "4 vs Will; 1d6 psychic damage. Repeat if even number."
$ws.Range("B14").Value = "Ranged 10; 4 vs Fortitude; 1d8 + 4 thunder damage."
$ws.Range("B17").Value = "1 or 2 creatures in close burst 3; 4 vs Will"
$ws.Range("B18").Value = "1d10 + 4 psychic damage, and push the target 4 squares. If"
$ws.Range("B21").Value = "Ranged 10; 4 vs Will; 6d6 + 4 modifier radiant damage."
$ws.Range("B23").Value = "a -4 penalty to attack rolls against me (save ends)."

# A33/B34: Skills line
Set-RichCell $ws.Range("A33") @(
    @{Text="Skills"; Bold=$true},
    @{Text=" 4 acro, 7 arc, 2 ath, 9 blu, 4 dip, 2 dun, 6 end, 2 heal,"; Bold=$false}
)
$ws.Range("B34").Value = "2 his, 7 ins, 4 int, 4 nat, 4 perc, 2 rel, 2 ste, 4 str, 4 thiev"

# Ability scores: modifiers bumped by +1
Set-RichCell $ws.Range("A35") @(
    @{Text="Str"; Bold=$true},
    @{Text=" 12 (+2)"; Bold=$false}
)
Set-RichCell $ws.Range("C35") @(
    @{Text="Dex"; Bold=$true},
    @{Text=" 16 (+4)"; Bold=$false}
)
Set-RichCell $ws.Range("D35") @(
    @{Text="Wis"; Bold=$true},
    @{Text=" 12 (+2)"; Bold=$false}
)
Set-RichCell $ws.Range("A36") @(
    @{Text="Con"; Bold=$true},
    @{Text=" 11 (+1)"; Bold=$false}
)
Set-RichCell $ws.Range("C36") @(
    @{Text="Int"; Bold=$true},
    @{Text=" 13 (+2)"; Bold=$false}
)
Set-RichCell $ws.Range("D36") @(
    @{Text="Cha"; Bold=$true},
    @{Text=" 16 (+4)"; Bold=$false}
)

# Saved view: scrolled down a bit further and selection moved to D37.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D37").Select()
